$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column C
$ws.Range("C1").Value = "Valor"

# Add values for column C (rows 3,4,5,7); rows 2 and 6 stay empty
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 5.01
$ws.Range("C5").Value = 7
$ws.Range("C7").Value = 9

# Apply currency cell style ("Moeda" / Currency) to column C data cells
$ws.Range("C1:C7").Style = "Currency"

# Adjust column widths
$ws.Columns.Item(2).ColumnWidth = 8.28515625
$ws.Columns.Item(3).ColumnWidth = 9.140625

# Update selection
$ws.Range("C6").Select()

# Update window view settings
$excel.ActiveWindow.WindowState = -4143
